# Generate Report for Handback
# Refreshes the timestamp values recorded in the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date for 6eb77c0c-...md
$wsOverview.Range("G2").Value = "2016-08-19 04:58:33"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime for 6eb77c0c-...md
$wsZhCn.Range("H2").Value = "2016-08-19 04:58:28"
$wsZhCn.Range("K2").Value = "2016-08-19 04:58:49"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime for 6eb77c0c-...md
$wsDeDe.Range("H2").Value = "2016-08-19 04:58:33"
$wsDeDe.Range("K2").Value = "2016-08-19 04:58:57"
